$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header (reuse the same header style as A1:E1)
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Convert A2:A7 timestamp strings into real datetime serial values
$dates = @(
    45687.48712708333,
    45687.52183310186,
    45687.52435972222,
    45687.52182962963,
    45687.52435625,
    45687.52965370371
)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]
    if ($i -eq 0) {
        $cell.NumberFormat = "yyyy-mm-dd h:mm:ss"
    }
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 6).Value = "Gra"
}
